$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualiza base de datos EC: intercambia los valores de "Valor Mora" entre
# el periodo 1911 (fila 16) y el periodo 1903 (fila 22).
$ws.Range("F16").Value = 33125
$ws.Range("F22").Value = 26500
